$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the text format to column A (rows 1-32) -- mirrors the
# upstream regeneration of the sheet's cell styles for this column.
$ws.Range("A1:A32").NumberFormat = "@"

# New prediction distance/score values for column B (rows 2-32),
# replacing the placeholder 1s written in the previous export.
$values = @(
  6422.5546610724014,
  4094.5597540970334,
  5840.9724929303447,
  4503.4900173480128,
  5527.3995076262818,
  5322.4579175880426,
  5497.2077339827356,
  9320.3721003090359,
  6723.2467083000265,
  9302.8779889828256,
  11112.987563362522,
  4657.974863404017,
  4486.4730165856708,
  4622.6759045851941,
  4869.0696839106131,
  4869.0696839106131,
  4869.0696839106131,
  4869.0696839106131,
  4292.3576440969773,
  4292.3576440969773,
  4869.0696839106131,
  4869.0696839106131,
  4869.0696839106131,
  4869.0696839106131,
  4869.0696839106131,
  4292.3576440969773,
  4869.0696839106131,
  4869.0696839106131,
  4292.3576440969773,
  9662.041362311038,
  5740.2400216305377
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value2 = $values[$i]
}
